$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.871.02"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.616.80"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  -0.83%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.60"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.38"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.840.73"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.614.18"
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.877.25"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.41"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.23"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.63"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.41"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.126.40"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.16"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.752.35"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.750"
$ws.Range("E43").Value = "  -5.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.03"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.95"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.410"
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.45"
$ws.Range("E51").Value = "  -1.90%  "
